# Natmi following Dr Hou advice
# Update ligand/receptor-expressing cell counts (from 1 to 3) and the
# resulting recalculated expression / specificity statistics for every
# data row (rows 2-7) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.141592333333333
$ws.Range("H2").Value = 3.424777
$ws.Range("I2").Value = 0.5720393989604073
$ws.Range("J2").Value = 0.5720393989604073
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 129.8307836666667
$ws.Range("N2").Value = 389.492351
$ws.Range("O2").Value = 0.4232533182703919
$ws.Range("P2").Value = 0.4232533182703919
$ws.Range("Q2").Value = 148.2138272645252
$ws.Range("R2").Value = 1333.924445380727
$ws.Range("S2").Value = 0.242117573791393
$ws.Range("T2").Value = 0.242117573791393

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.141592333333333
$ws.Range("H3").Value = 3.424777
$ws.Range("I3").Value = 0.5720393989604073
$ws.Range("J3").Value = 0.5720393989604073
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 118.3062056666667
$ws.Range("N3").Value = 354.918617
$ws.Range("O3").Value = 0.3856827533981234
$ws.Range("P3").Value = 0.3856827533981234
$ws.Range("Q3").Value = 135.0574573748232
$ws.Range("R3").Value = 1215.517116373409
$ws.Range("S3").Value = 0.2206257304432575
$ws.Range("T3").Value = 0.2206257304432575

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.141592333333333
$ws.Range("H4").Value = 3.424777
$ws.Range("I4").Value = 0.5720393989604073
$ws.Range("J4").Value = 0.5720393989604073
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 58.60787966666667
$ws.Range("N4").Value = 175.823639
$ws.Range("O4").Value = 0.1910639283314847
$ws.Range("P4").Value = 0.1910639283314848
$ws.Range("Q4").Value = 66.90630610038923
$ws.Range("R4").Value = 602.156754903503
$ws.Range("S4").Value = 0.1092960947257569
$ws.Range("T4").Value = 0.1092960947257569

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.8540610000000001
$ws.Range("H5").Value = 2.562183
$ws.Range("I5").Value = 0.4279606010395928
$ws.Range("J5").Value = 0.4279606010395928
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 129.8307836666667
$ws.Range("N5").Value = 389.492351
$ws.Range("O5").Value = 0.4232533182703919
$ws.Range("P5").Value = 0.4232533182703919
$ws.Range("Q5").Value = 110.883408929137
$ws.Range("R5").Value = 997.950680362233
$ws.Range("S5").Value = 0.181135744478999
$ws.Range("T5").Value = 0.181135744478999

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.8540610000000001
$ws.Range("H6").Value = 2.562183
$ws.Range("I6").Value = 0.4279606010395928
$ws.Range("J6").Value = 0.4279606010395928
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 118.3062056666667
$ws.Range("N6").Value = 354.918617
$ws.Range("O6").Value = 0.3856827533981234
$ws.Range("P6").Value = 0.3856827533981234
$ws.Range("Q6").Value = 101.040716317879
$ws.Range("R6").Value = 909.366446860911
$ws.Range("S6").Value = 0.1650570229548659
$ws.Range("T6").Value = 0.1650570229548659

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.8540610000000001
$ws.Range("H7").Value = 2.562183
$ws.Range("I7").Value = 0.4279606010395928
$ws.Range("J7").Value = 0.4279606010395928
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 58.60787966666667
$ws.Range("N7").Value = 175.823639
$ws.Range("O7").Value = 0.1910639283314847
$ws.Range("P7").Value = 0.1910639283314848
$ws.Range("Q7").Value = 50.05470431599301
$ws.Range("R7").Value = 450.492338843937
$ws.Range("S7").Value = 0.08176783360572788
$ws.Range("T7").Value = 0.0817678336057279
